$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.567.76"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "1.472.24"
$ws.Range("E3").Value = "  +2.09%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9582"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3563"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3070"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.090"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06637"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.458"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.177"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9582"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.97%  "

$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").Value = "1.471.47"
$ws.Range("E18").Value = "  +2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05959"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.488"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.273"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "

$ws.Range("D25").Value = "20.557.89"
$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.29%  "

$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("D29").Value = "1.630.79"
$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.860"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.933"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("E33").Value = "  +3.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7946"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.443"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05735"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.711"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02029"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9590"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.42%  "

$ws.Range("E41").Value = "  +1.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.279"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5248"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.509"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5183"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.801"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06436"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9911"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.68%  "
